$wb = $excel.ActiveWorkbook

# Grab the two worksheets by their current names.
$wsHotel  = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# Reorder the tabs so that "review_info" is now the first sheet and
# "hotel_info" is the second sheet.
$wsReview.Move($wsHotel)

# Worksheet handles obtained before the Move can now resolve to a
# different sheet (they track tab position), so re-fetch a fresh
# reference to "hotel_info" by name before editing it further.
$wsHotel = $wb.Worksheets.Item("hotel_info")

# Insert a new "State" column into hotel_info, between "Hotel_Name"
# (column B) and "City" (column C), and populate it.
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"
